# Generate Report for Handback
#
# The handoff files previously recorded in columns A (Source/md) and D
# (Latest Handoff File/xlf) have now been handed back. This records the
# same files as the "Latest Target File" (F) / "Latest Handback File" (G)
# pair, flips the Status message, and stamps the handback timestamp into
# "Latest Handback DateTime" (H) for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# OLE_COLOR (0x00BBGGRR) equivalent of the existing HyperLink font color
# FF6495ED (R=0x64 G=0x95 B=0xED) used by columns A/D so the new cells
# visually match the existing handoff-file hyperlinks.
$hyperlinkColor = 15570276

function Style-AsHyperlink($cell) {
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = $true
    $cell.Font.Color = $hyperlinkColor
}

function Set-HandbackColumns($ws, $mdUrl2, $xlfUrl2, $mdUrl3, $xlfUrl3, $handbackDateTime) {
    # Status: handoff is complete, now in sync with en-US source.
    $ws.Range("C2").Value2 = "Handed back: in sync with en-US"
    $ws.Range("C3").Value2 = "Handed back: in sync with en-US"

    # Row 2 ("2e2121cf-...") target/handback file columns mirror A2/D2.
    $mdName2 = $ws.Range("A2").Value2
    $xlfName2 = $ws.Range("D2").Value2

    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl2, "", "", $mdName2)
    Style-AsHyperlink($ws.Range("F2"))

    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl2, "", "", $xlfName2)
    Style-AsHyperlink($ws.Range("G2"))

    # Row 3 ("e0830705-...") target/handback file columns mirror A3/D3.
    $mdName3 = $ws.Range("A3").Value2
    $xlfName3 = $ws.Range("D3").Value2

    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl3, "", "", $mdName3)
    Style-AsHyperlink($ws.Range("F3"))

    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl3, "", "", $xlfName3)
    Style-AsHyperlink($ws.Range("G3"))

    # Latest Handback DateTime.
    $ws.Range("H2").Value2 = $handbackDateTime
    $ws.Range("H3").Value2 = $handbackDateTime
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackColumns $wsZhCn `
    "https://github.com/OpenLocalizationTest/oltest/blob/a8566d6605bc32d080b4135178549f507a29f532/e2e/2e2121cf-0973-468a-abc6-6648e40db18e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1844a5ea1413ba82c87cbdcda73fb27c7cc7d38b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/2e2121cf-0973-468a-abc6-6648e40db18e.b4fcd7cc55a65664f18afad341b35477176fea58.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a8566d6605bc32d080b4135178549f507a29f532/e2e/e0830705-ca48-4932-a975-71fa95d84373.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1844a5ea1413ba82c87cbdcda73fb27c7cc7d38b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/e0830705-ca48-4932-a975-71fa95d84373.ce6f8af22a9a751da26d6ac0dd9598add2200d1b.zh-cn.xlf" `
    "2016-03-24 07:47:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackColumns $wsDeDe `
    "https://github.com/OpenLocalizationTest/oltest/blob/a8566d6605bc32d080b4135178549f507a29f532/e2e/2e2121cf-0973-468a-abc6-6648e40db18e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/158f324dfbc86522f08b56552b737c90dd36e96d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/2e2121cf-0973-468a-abc6-6648e40db18e.b4fcd7cc55a65664f18afad341b35477176fea58.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a8566d6605bc32d080b4135178549f507a29f532/e2e/e0830705-ca48-4932-a975-71fa95d84373.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/158f324dfbc86522f08b56552b737c90dd36e96d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/e0830705-ca48-4932-a975-71fa95d84373.ce6f8af22a9a751da26d6ac0dd9598add2200d1b.de-de.xlf" `
    "2016-03-24 07:47:27"
